$wb = $excel.ActiveWorkbook

# Rename the "Include from Reasons why cert" worksheet tab to "Include #0"
$wsInclude = $wb.Worksheets.Item("Include from Reasons why cert")
$wsInclude.Name = "Include #0"

# --- Metadata sheet updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Update URL value (pythia -> cicada)
$wsMeta.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/preferred-allowed-reason"

# Update Date value
$wsMeta.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# Insert a new row after row 10 (Contact) for Jurisdiction
$wsMeta.Rows.Item(11).Insert()
$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""

# Copy style from row 10 (existing data row) into new row 11 so formatting matches
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122) | Out-Null

# --- Include sheet updates ---
$wsInclude.Range("B9").Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/PreferredAllowedReason"
